$wb = $excel.ActiveWorkbook

# Update the "想去人数" (want-to-go count) figures for the three events that
# were refreshed in this data snapshot. The same three rows appear both on
# the "展览" sheet and on the "全部类型" sheet, so both need updating.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 261
    $ws.Range("F3").Value = 84
    $ws.Range("F4").Value = 919
}
